$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")
$ws.Activate()

# Update the repaymentstrategy value (row 17) from "RBI (India)" to the
# new periodic/upfront scenario value.
$ws.Range("B17").Value = "Overdue/Due Fee/Int,Principal"
$ws.Range("B17").Select()
